$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 56: "Week 7" header (plain / default style) ---
$ws.Range("A56").Value = "Week 7"

# --- New row 59/60 block: orange-highlighted "Integrate User Model..." pair, reusing the same
#     font treatment already used at B52/B54 (fontId 4 -- plain font, explicitly applied) ---
$ws.Range("B52").Copy()
$ws.Range("A59:B60").PasteSpecial(-4122)
$ws.Range("A59:B60").Interior.Color = 49407

$ws.Range("A59").Value = "Integrate User Model with Existing Apps"
$ws.Range("B59").Value = "TimeTracker"
$ws.Range("B60").Value = "Calendar"

# --- Row 57/58: existing "Manage permissions..." pair gets a red highlight ---
$ws.Range("A57").Interior.Color = 255
$ws.Range("A58").Interior.Color = 255

# --- New row 61: "Settings Panel", orange highlight (plain font, reuses the orange fill) ---
$ws.Range("A61").Interior.Color = 49407
$ws.Range("A61").Value = "Settings Panel"

# --- New row 62: " Hierarchy With Existing User Levels", red highlight, reusing the red fill/style ---
$ws.Range("A62").Interior.Color = 255
$ws.Range("A62").Value = " Hierarchy With Existing User Levels"

# --- Column C notes (people assigned) ---
$ws.Range("C61").Value = "Sid + Shivam"
$ws.Range("C62").Value = "Lohit + Arnav"
$ws.Range("C59").Value = "Arny + Shivam"
$ws.Range("C57").Value = "Lohit (due to refusal to contact Sid)"

# --- Center-align & merge the C59:C60 note cell ---
$ws.Range("C59:C60").HorizontalAlignment = -4108
$ws.Range("C59:C60").Merge()

# --- Scroll / selection matches the author's final view ---
$ws.Application.ActiveWindow.ScrollRow = 46
$ws.Range("C64").Select()
